$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Update the Date value cell (row 8, column B) to the new timestamp
$meta.Range("B8").Value = "2022-05-13T16:07:55-04:00"

# Fill in the Copyright value cell (row 15, column B) with the LOINC/SNOMED copyright text
$meta.Range("B15").Value = "This material contains content from LOINC (http://loinc.org). LOINC is copyright © 1995-2020, Regenstrief Institute, Inc. and the Logical Observation Identifiers Names and Codes (LOINC) Committee and is available at no cost under the license at http://loinc.org/license. LOINC® is a registered United States trademark of Regenstrief Institute, Inc. This value set includes content from SNOMED CT, which is copyright © 2002+ International Health Terminology Standards Development Organisation (IHTSDO), and distributed by agreement between IHTSDO and HL7. Implementer use of SNOMED CT is not covered by this agreement"
